# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 7464
    3  = 7529
    4  = 106
    8  = 121
    9  = 99
    10 = 141
    12 = 102
    13 = 680
    14 = 610
    16 = 38
    18 = 129
    19 = 87
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
